$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---- Header row (row 1): add proper column labels ----
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Column I ("date") must stay literal text, not get auto-parsed into a date serial.
# Temporarily force a text format so typing "2013-12-31" isn't reinterpreted,
# then clear the format again so the cell ends up with the plain/default style
# (matching the rest of the row) while keeping the text value.
$ws.Range("I2:I5").NumberFormat = "@"

# ---- Data rows 2-5: append the new metadata columns G-M ----
# (columns B-F already hold the correct bank/deposit_type/currency/owner/total
#  values; row 5 col F is retyped from text to a real number below)
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 7).Value  = "deposit"      # property_category
    $ws.Cells.Item($r, 8).Value  = "normal"       # category
    $ws.Cells.Item($r, 9).Value  = "2013-12-31"   # date
    $ws.Cells.Item($r, 10).Value = "姚文智"        # legislator_name
    $ws.Cells.Item($r, 11).Value = 1745           # legislator_id
    $ws.Cells.Item($r, 12).Value = "tmpc2191"     # source_file
    $ws.Cells.Item($r, 13).Value = $r + 46        # index (48,49,50,51)
}

# Row 5's "total" (F5) was stored as text "2998776" - normalize to a number
$ws.Range("F5").Value = 2998776

# Drop the temporary "@" text format now that the literal text is locked in,
# so these cells fall back to the same plain/default style as their neighbours
$ws.Range("I2:I5").ClearFormats()

# ---- Formatting: match the existing header style (bold + border + centered)
#      on the new header cells G1:M1. The new data cells G2:M5 keep the
#      plain/default formatting already shared by the rest of the table. ----
$newHeaderRange = $ws.Range("G1:M1")
$newHeaderRange.Font.Bold = $true
$newHeaderRange.Borders.LineStyle = 1
$newHeaderRange.HorizontalAlignment = -4108
$newHeaderRange.VerticalAlignment = -4160
